# Auto update Excel log — append newly-logged sensor rows (2026-02-06, 10:18-10:19)
# to the PIR, Humidity and Temperature sheets.

$wb = $excel.ActiveWorkbook

# --- PIR: append rows 460-472 ---
$ws = $wb.Worksheets.Item("PIR")
$newRows = @(
    ,("'2026-02-06", "10:18:07", "10:00", "Bathroom", "No Motion", "Inactive")
    ,("'2026-02-06", "10:18:10", "10:00", "Bathroom", "No Motion", "Inactive")
    ,("'2026-02-06", "10:18:14", "10:00", "Bathroom", "No Motion", "Inactive")
    ,("'2026-02-06", "10:18:19", "10:00", "Bathroom", "No Motion", "Inactive")
    ,("'2026-02-06", "10:18:24", "10:00", "Bathroom", "No Motion", "Inactive")
    ,("'2026-02-06", "10:18:29", "10:00", "Bathroom", "No Motion", "Inactive")
    ,("'2026-02-06", "10:18:34", "10:00", "Bathroom", "No Motion", "Inactive")
    ,("'2026-02-06", "10:18:39", "10:00", "Bathroom", "No Motion", "Inactive")
    ,("'2026-02-06", "10:18:44", "10:00", "Bathroom", "No Motion", "Inactive")
    ,("'2026-02-06", "10:18:49", "10:00", "Bathroom", "No Motion", "Inactive")
    ,("'2026-02-06", "10:18:54", "10:00", "Bathroom", "No Motion", "Inactive")
    ,("'2026-02-06", "10:18:59", "10:00", "Bathroom", "No Motion", "Inactive")
    ,("'2026-02-06", "10:19:04", "10:00", "Bathroom", "No Motion", "Inactive")
)
$r = 460
foreach ($row in $newRows) {
    for ($c = 0; $c -lt 6; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

# --- Humidity: append rows 315-326 ---
$ws = $wb.Worksheets.Item("Humidity")
$newRows = @(
    ,("'2026-02-06", "10:18:08", "10:00", "Bathroom", "'68.8%", "Active")
    ,("'2026-02-06", "10:18:11", "10:00", "Bathroom", "'68.8%", "Active")
    ,("'2026-02-06", "10:18:15", "10:00", "Bathroom", "'68.7%", "Active")
    ,("'2026-02-06", "10:18:20", "10:00", "Bathroom", "'68.9%", "Active")
    ,("'2026-02-06", "10:18:25", "10:00", "Bathroom", "'68.9%", "Active")
    ,("'2026-02-06", "10:18:30", "10:00", "Bathroom", "'69.0%", "Active")
    ,("'2026-02-06", "10:18:35", "10:00", "Bathroom", "'69.1%", "Active")
    ,("'2026-02-06", "10:18:40", "10:00", "Bathroom", "'69.0%", "Active")
    ,("'2026-02-06", "10:18:45", "10:00", "Bathroom", "'68.9%", "Active")
    ,("'2026-02-06", "10:18:50", "10:00", "Bathroom", "'68.8%", "Active")
    ,("'2026-02-06", "10:18:55", "10:00", "Bathroom", "'68.8%", "Active")
    ,("'2026-02-06", "10:19:06", "10:00", "Bathroom", "'68.5%", "Active")
)
$r = 315
foreach ($row in $newRows) {
    for ($c = 0; $c -lt 6; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

# --- Temperature: append rows 315-326 ---
$ws = $wb.Worksheets.Item("Temperature")
$newRows = @(
    ,("'2026-02-06", "10:18:09", "10:00", "Bathroom", "28.2C", "Active")
    ,("'2026-02-06", "10:18:12", "10:00", "Bathroom", "28.2C", "Active")
    ,("'2026-02-06", "10:18:16", "10:00", "Bathroom", "28.1C", "Active")
    ,("'2026-02-06", "10:18:21", "10:00", "Bathroom", "28.2C", "Active")
    ,("'2026-02-06", "10:18:26", "10:00", "Bathroom", "28.2C", "Active")
    ,("'2026-02-06", "10:18:31", "10:00", "Bathroom", "28.1C", "Active")
    ,("'2026-02-06", "10:18:36", "10:00", "Bathroom", "28.2C", "Active")
    ,("'2026-02-06", "10:18:41", "10:00", "Bathroom", "28.2C", "Active")
    ,("'2026-02-06", "10:18:46", "10:00", "Bathroom", "28.2C", "Active")
    ,("'2026-02-06", "10:18:51", "10:00", "Bathroom", "28.3C", "Active")
    ,("'2026-02-06", "10:18:56", "10:00", "Bathroom", "28.3C", "Active")
    ,("'2026-02-06", "10:19:07", "10:00", "Bathroom", "28.2C", "Active")
)
$r = 315
foreach ($row in $newRows) {
    for ($c = 0; $c -lt 6; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}
